$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 92
$ws.Range("H92").Value = 755.7857
$ws.Range("I92").Value = 327.6
$ws.Range("K92").Value = 327.6
$ws.Range("M92").Value = 920.4
# Row 132
$ws.Range("H132").Value = 4493.4287
$ws.Range("I132").Value = 2897.2666
$ws.Range("K132").Value = 8691.799800000001
$ws.Range("M132").Value = -6161.799800000001
# Row 137
$ws.Range("H137").Value = 3622.7144
$ws.Range("I137").Value = 2976.5
$ws.Range("K137").Value = 8929.5
$ws.Range("M137").Value = -6379.5

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 2054.919
$ws.Range("I2").Value = 2162.0344
$ws.Range("J2").Value = 1666.625
$ws.Range("K2").Value = 2162.0344
$ws.Range("L2").Value = 1666.625
$ws.Range("M2").Value = -2049.0344
$ws.Range("N2").Value = -1892.625
# Row 32
$ws.Range("H32").Value = 4563.904
$ws.Range("I32").Value = 3965.7754
$ws.Range("K32").Value = 3965.7754
$ws.Range("M32").Value = -3678.7754
# Row 74
$ws.Range("H74").Value = 4053.1875
$ws.Range("I74").Value = 3688.5833
$ws.Range("K74").Value = 3688.5833
$ws.Range("M74").Value = -2814.5833
# Row 77
$ws.Range("H77").Value = 4053.1875
$ws.Range("I77").Value = 3688.5833
$ws.Range("K77").Value = 18442.9165
$ws.Range("M77").Value = -14074.9165
# Row 97
$ws.Range("H97").Value = 1410.92
$ws.Range("I97").Value = 1256.2778
$ws.Range("J97").Value = 1808.5714
$ws.Range("K97").Value = 1256.2778
$ws.Range("L97").Value = 1808.5714
$ws.Range("M97").Value = -760.2778000000001
$ws.Range("N97").Value = -2800.5714
# Row 102
$ws.Range("H102").Value = 1834.2354
$ws.Range("I102").Value = 1307.2727
$ws.Range("K102").Value = 1307.2727
$ws.Range("M102").Value = 314.7273
# Row 116
$ws.Range("H116").Value = 2054.919
$ws.Range("I116").Value = 2162.0344
$ws.Range("J116").Value = 1666.625
$ws.Range("K116").Value = 2162.0344
$ws.Range("L116").Value = 1666.625
$ws.Range("M116").Value = 131.9656
$ws.Range("N116").Value = -6254.625
# Row 122
$ws.Range("H122").Value = 3400.9092
$ws.Range("I122").Value = 3547.15
$ws.Range("K122").Value = 10641.45
$ws.Range("M122").Value = -8191.450000000001
# Row 132
$ws.Range("H132").Value = 2504044.2
$ws.Range("I132").Value = 4105.2583
$ws.Range("K132").Value = 12315.7749
$ws.Range("M132").Value = -9785.7749

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 2054.919
$ws.Range("I3").Value = 2162.0344
$ws.Range("J3").Value = 1666.625
$ws.Range("K3").Value = 2162.0344
$ws.Range("L3").Value = 1666.625
$ws.Range("M3").Value = -2048.0344
$ws.Range("N3").Value = -1894.625
# Row 20
$ws.Range("H20").Value = 11303.111
$ws.Range("J20").Value = 765.6667
$ws.Range("L20").Value = 765.6667
$ws.Range("N20").Value = -1259.6667
# Row 86
$ws.Range("H86").Value = 3177.6052
$ws.Range("I86").Value = 1451
$ws.Range("K86").Value = 1451
$ws.Range("M86").Value = -328
# Row 89
$ws.Range("H89").Value = 3177.6052
$ws.Range("I89").Value = 1451
$ws.Range("K89").Value = 7255
$ws.Range("M89").Value = -1639
# Row 107
$ws.Range("H107").Value = 4916.615
$ws.Range("I107").Value = 5519.6
$ws.Range("K107").Value = 5519.6
$ws.Range("M107").Value = -3599.6

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 3848474
$ws.Range("I16").Value = 4350131.5
$ws.Range("K16").Value = 4350131.5
$ws.Range("M16").Value = -4349844.5
# Row 107
$ws.Range("H107").Value = 1232.2084
$ws.Range("I107").Value = 880.13635
$ws.Range("K107").Value = 880.13635
$ws.Range("M107").Value = 1039.86365
# Row 113
$ws.Range("H113").Value = 3848474
$ws.Range("I113").Value = 4350131.5
$ws.Range("K113").Value = 4350131.5
$ws.Range("M113").Value = -4347961.5
# Row 134
$ws.Range("H134").Value = 2106.1
$ws.Range("I134").Value = 1978
$ws.Range("J134").Value = 2191.5
$ws.Range("K134").Value = 5934
$ws.Range("L134").Value = 6574.5
$ws.Range("M134").Value = -3399
$ws.Range("N134").Value = -11644.5

$ws = $wb.Worksheets.Item("CUL")
# Row 33
$ws.Range("H33").Value = 2482272.5
$ws.Range("I33").Value = 365
$ws.Range("J33").Value = 4033464.5
$ws.Range("K33").Value = 2190
$ws.Range("L33").Value = 24200787
$ws.Range("M33").Value = -1907
$ws.Range("N33").Value = -24201353

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 56
$ws.Range("I2").Value = 60.5
$ws.Range("J2").Value = 50
$ws.Range("K2").Value = 60.5
$ws.Range("L2").Value = 50
$ws.Range("M2").Value = 52.5
$ws.Range("N2").Value = -276
# Row 70
$ws.Range("H70").Value = 9358.421
$ws.Range("I70").Value = 9067.091
$ws.Range("K70").Value = 9067.091
$ws.Range("M70").Value = -8797.091
# Row 73
$ws.Range("H73").Value = 9358.421
$ws.Range("I73").Value = 9067.091
$ws.Range("K73").Value = 9067.091
$ws.Range("M73").Value = -8131.091
# Row 102
$ws.Range("H102").Value = 2502.4688
$ws.Range("I102").Value = 2468.2144
$ws.Range("J102").Value = 2742.25
$ws.Range("K102").Value = 2468.2144
$ws.Range("L102").Value = 2742.25
$ws.Range("M102").Value = -846.2143999999998
$ws.Range("N102").Value = -5986.25
# Row 132
$ws.Range("H132").Value = 5002083.5
$ws.Range("I132").Value = 1844.4667
$ws.Range("K132").Value = 5533.4001
$ws.Range("M132").Value = -3003.4001

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 4269.7144
$ws.Range("I7").Value = 3939.7144
$ws.Range("J7").Value = 4929.7144
$ws.Range("K7").Value = 3939.7144
$ws.Range("L7").Value = 4929.7144
$ws.Range("M7").Value = -3827.7144
$ws.Range("N7").Value = -5153.7144
# Row 40
$ws.Range("H40").Value = 3550.3333
$ws.Range("I40").Value = 3308.9443
$ws.Range("J40").Value = 4998.6665
$ws.Range("K40").Value = 3308.9443
$ws.Range("L40").Value = 4998.6665
$ws.Range("M40").Value = -3172.9443
$ws.Range("N40").Value = -5270.6665
# Row 61
$ws.Range("H61").Value = 3408.3333
$ws.Range("I61").Value = 2613
$ws.Range("K61").Value = 2613
$ws.Range("M61").Value = -2411
# Row 113
$ws.Range("H113").Value = 3408.3333
$ws.Range("I113").Value = 2613
$ws.Range("K113").Value = 2613
$ws.Range("M113").Value = -443
# Row 126
$ws.Range("H126").Value = 4269.7144
$ws.Range("I126").Value = 3939.7144
$ws.Range("J126").Value = 4929.7144
$ws.Range("K126").Value = 11819.1432
$ws.Range("L126").Value = 14789.1432
$ws.Range("M126").Value = -9349.143199999999
$ws.Range("N126").Value = -19729.1432
# Row 136
$ws.Range("H136").Value = 2888.44
$ws.Range("I136").Value = 2874
$ws.Range("K136").Value = 8622
$ws.Range("M136").Value = -6072

$ws = $wb.Worksheets.Item("WVR")
# Row 4
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").ClearContents()
# Row 122
$ws.Range("H122").Value = 2134.55
$ws.Range("I122").Value = 2214.3076
$ws.Range("J122").Value = 1986.4286
$ws.Range("K122").Value = 6642.9228
$ws.Range("L122").Value = 5959.2858
$ws.Range("M122").Value = -4192.9228
$ws.Range("N122").Value = -10859.2858
# Row 126
$ws.Range("H126").Value = 3400.4
$ws.Range("I126").Value = 3646.6155
$ws.Range("K126").Value = 10939.8465
$ws.Range("M126").Value = -8469.8465

Write-Output "Applied Ragnarok_Profits scheduled update"